$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 (treatment): update recomputed meta-analysis values ---
$ws1.Range("B3").Value = 0.5002105923365068
$ws1.Range("C3").Value = 0.7786591166611616
$ws1.Range("D3").Value = 1.157264633490854

$ws1.Range("B4").Value = 0.1205578087620536
$ws1.Range("C4").Value = 0.455431433919275
$ws1.Range("D4").Value = 1.008734676105219

$ws1.Range("B5").Value = 0.3582556579487085
$ws1.Range("C5").Value = 0.6963168204826793
$ws1.Range("D5").Value = 1.036296091716069

$ws1.Range("B7").Value = 17.10110608132888
$ws1.Range("C7").Value = 25.87636470345233
$ws1.Range("D7").Value = 37.61825124372285

$ws1.Range("B8").Value = 0.114559495474643
$ws1.Range("C8").Value = 0.4115696663083406
$ws1.Range("D8").Value = 0.8950197316576815

$ws1.Range("B9").Value = 0.3486282232034319
$ws1.Range("C9").Value = 0.6607979328218138
$ws1.Range("D9").Value = 0.9744583947394193

$ws1.Range("B11").Value = 8.056419274879348
$ws1.Range("C11").Value = 10.90026959821771
$ws1.Range("D11").Value = 14.40380497647034

$ws1.Range("B12").Value = 0.05255909735591642
$ws1.Range("C12").Value = 0.2027227971071784
$ws1.Range("D12").Value = 0.4523595204742347

$ws1.Range("B13").Value = 0.2367176417544625
$ws1.Range("C13").Value = 0.4648983143391469
$ws1.Range("D13").Value = 0.6944620917049884

# New: Speed meta analysis section
$ws1.Range("A14").Value = "Speed meta analysis"

$ws1.Range("A15").Value = "mean (km/day)"
$ws1.Range("B15").Value = 2.97132081279445
$ws1.Range("C15").Value = 3.753048655222063
$ws1.Range("D15").Value = 4.676696035531339

$ws1.Range("A16").Value = "CoV² (RVAR)"
$ws1.Range("B16").Value = 0.03301988768423263
$ws1.Range("C16").Value = 0.1193087915423664
$ws1.Range("D16").Value = 0.2600045075943561

$ws1.Range("A17").Value = "CoV  (RSTD)"
$ws1.Range("B17").Value = 0.187205985989238
$ws1.Range("C17").Value = 0.355850974056521
$ws1.Range("D17").Value = 0.5253181112768256

# --- Sheet2 (control): update recomputed meta-analysis values ---
$ws2.Range("B3").Value = 0.5178416068095788
$ws2.Range("C3").Value = 0.9296259217781964
$ws2.Range("D3").Value = 1.545114629883028

$ws2.Range("B4").Value = 0.07253888170390142
$ws2.Range("C4").Value = 0.4774612343813659
$ws2.Range("D4").Value = 1.25591472488125

$ws2.Range("B5").Value = 0.2833810506161001
$ws2.Range("C5").Value = 0.7270335229971665
$ws2.Range("D5").Value = 1.17914022275735

$ws2.Range("B7").Value = 22.84308894697649
$ws2.Range("C7").Value = 35.34294070428853
$ws2.Range("D7").Value = 52.31257022832333

$ws2.Range("B8").Value = 0.04750194107559804
$ws2.Range("C8").Value = 0.2723737122129692
$ws2.Range("D8").Value = 0.6896451688430465

$ws2.Range("B9").Value = 0.2281988506563955
$ws2.Range("C9").Value = 0.5464373410218864
$ws2.Range("D9").Value = 0.8695020765546931

# New: Speed meta analysis section
$ws2.Range("A14").Value = "Speed meta analysis"

$ws2.Range("A15").Value = "mean (km/day)"
$ws2.Range("B15").Value = 3.5279196661446
$ws2.Range("C15").Value = 4.237933798281795
$ws2.Range("D15").Value = 5.045164642902208

$ws2.Range("A16").Value = "CoV² (RVAR)"
$ws2.Range("B16").Value = 0.005267892284647601
$ws2.Range("C16").Value = 0.03804643118953414
$ws2.Range("D16").Value = 0.1025332284868073

$ws2.Range("A17").Value = "CoV  (RSTD)"
$ws2.Range("B17").Value = 0.07661869813981735
$ws2.Range("C17").Value = 0.2059080321739013
$ws2.Range("D17").Value = 0.3380248599474686
